$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.216.47"
$ws.Range("E2").Value = "  -0.36%  "
$ws.Range("D3").Value = "1.928.11"
$ws.Range("E3").Value = "  -0.31%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'248.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.14%  "
$ws.Range("D6").Value = "'0.7144"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.21%  "
$ws.Range("D7").Value = "'1.000"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "'0.3211"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.37%  "
$ws.Range("D9").Value = "'27.45"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.50%  "
$ws.Range("D10").Value = "'0.07098"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.42%  "
$ws.Range("D11").Value = "'0.7926"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.64%  "
$ws.Range("D12").Value = "'0.07951"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.57%  "
$ws.Range("D13").Value = "1.927.08"
$ws.Range("D14").Value = "'5.381"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.58%  "
$ws.Range("D15").Value = "'94.87"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.03%  "
$ws.Range("D16").Value = "'14.68"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.97%  "
$ws.Range("D17").Value = "30.217.91"
$ws.Range("E17").Value = "  -0.32%  "
$ws.Range("D18").Value = "'256.10"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.31%  "
$ws.Range("D19").Value = "'0.000008054"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.79%  "
$ws.Range("D20").Value = "'5.762"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.10%  "
$ws.Range("D21").Value = "2.181.23"
$ws.Range("E21").Value = "  -0.17%  "
$ws.Range("D22").Value = "'0.9998"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("D23").Value = "'0.9998"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").Value = "'6.833"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.57%  "
$ws.Range("D25").Value = "'9.530"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.98%  "
$ws.Range("D26").Value = "'165.95"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.03%  "
$ws.Range("D27").Value = "'19.08"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.33%  "
$ws.Range("D28").Value = "'2.272"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.44%  "
$ws.Range("D29").Value = "'0.1268"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.39%  "
$ws.Range("D30").Value = "'1.356"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.33%  "
$ws.Range("D31").Value = "'1.529"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.98%  "
$ws.Range("D32").Value = "'4.395"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.16%  "
$ws.Range("D33").Value = "'4.130"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.40%  "
$ws.Range("D34").Value = "'0.05153"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.87%  "
$ws.Range("D35").Value = "'1.270"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.84%  "
$ws.Range("D36").Value = "'0.7462"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.62%  "
$ws.Range("D37").Value = "'2.759"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.32%  "
$ws.Range("D38").Value = "'0.01956"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.74%  "
$ws.Range("D39").Value = "'2.797"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.36%  "
$ws.Range("D40").Value = "'77.70"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.61%  "
$ws.Range("D41").Value = "'6.361"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.02%  "
$ws.Range("D42").Value = "'0.4501"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.73%  "
$ws.Range("D43").Value = "'1.984"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.72%  "
$ws.Range("D44").Value = "'0.8439"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.97%  "
$ws.Range("D45").Value = "'0.9995"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.11%  "
$ws.Range("D46").Value = "'100.61"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.47%  "
$ws.Range("D47").Value = "'9.723"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.57%  "
$ws.Range("D48").Value = "'7.435"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.58%  "
$ws.Range("D49").Value = "'36.41"
$ws.Range("D49").Style = "Normal"

# Row 50: SynthetixNetwork -> Cronos
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.06101"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.41%  "

# Row 51: Cronos -> Decentraland
$ws.Range("B51").Value = "Decentraland"
$ws.Range("C51").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D51").Value = "'0.4201"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.16%  "

